# Generate Report for Handback
#
# The handback CI job re-ran and produced a new pair of generated file
# records (new GUID-based file names / content hashes) together with
# refreshed timestamps. This replaces the two rows reported on the
# Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "ac2f5080-7f2a-49cf-be49-8ed8dacc307e"
$newGuid1 = "5a57b781-858c-4266-b122-ad0635dcfa74"
$oldGuid2 = "afb46ee5-6896-4257-a56d-04be8f8c5f92"
$newGuid2 = "ffffdb97a82b-372e-4208-ad5f-b8e8181f25c3"

$newHash = "92d3c8e7f0eaf37956d44926cd72a9d6ac625084"

$newOverviewDate = "2016-08-23 00:59:04"

$newZhHandoffDate = "2016-08-23 00:58:56"
$newZhHandbackDate = "2016-08-23 00:59:27"

$newDeHandoffDate = "2016-08-23 00:59:04"
$newDeHandbackDate = "2016-08-23 00:59:33"

# ---------- Overview sheet ----------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid1.md"
$wsOverview.Range("G2").Value = $newOverviewDate

$wsOverview.Range("A3").Value = "$newGuid2.md"
$wsOverview.Range("G3").Value = $newOverviewDate

# Hyperlinks.Delete() on this engine clears the whole sheet collection, so
# rebuild all of them (preserving their original target URLs - only the
# displayed text changes) in the original order so relationship ids land
# back on rId2/rId3/...
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid1.md", "", "", "e2e\$newGuid1.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid2.md", "", "", "e2e\$newGuid2.md") | Out-Null

# ---------- zh-cn sheet ----------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("I2").Value = "$newGuid1.md"
$wsZh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("J2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = $newZhHandoffDate
$wsZh.Range("K2").Value = $newZhHandbackDate

$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Range("I3").Value = "$newGuid2.md"
$wsZh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("J3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$wsZh.Range("H3").Value = $newZhHandoffDate
$wsZh.Range("K3").Value = $newZhHandbackDate

$wsZh.Range("A1").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/24538a3aacb9fb6f651ce64f8da5bdd49ced9576/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/24538a3aacb9fb6f651ce64f8da5bdd49ced9576/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null

# ---------- de-de sheet ----------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("I2").Value = "$newGuid1.md"
$wsDe.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("J2").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = $newDeHandoffDate
$wsDe.Range("K2").Value = $newDeHandbackDate

$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Range("I3").Value = "$newGuid2.md"
$wsDe.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("J3").Value = "$newGuid1.$newHash.de-de.xlf"
$wsDe.Range("H3").Value = $newDeHandoffDate
$wsDe.Range("K3").Value = $newDeHandbackDate

$wsDe.Range("A1").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d016b5d3e77bcf6eb2ac6a291db92e3f6455bb2a/e2e/$oldGuid1.md", "", "", "$newGuid1.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d016b5d3e77bcf6eb2ac6a291db92e3f6455bb2a/e2e/$oldGuid2.md", "", "", "$newGuid2.md") | Out-Null
